$d = $word.ActiveDocument

# Remove the parenthetical "(r > .7) " that followed
# "...correlate highly with one another " and preceded "then we will average them..."
$d.Content.Find.Execute("(r > .7) ", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "", 2)
